# Update "想去人数" (interest-count) figures on the "展览" sheet and on the
# aggregated "全部类型" sheet to match the freshly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 242
$wsExpo.Range("F6").Value  = 9995
$wsExpo.Range("F8").Value  = 331
$wsExpo.Range("F9").Value  = 1244
$wsExpo.Range("F10").Value = 5340
$wsExpo.Range("F14").Value = 121
$wsExpo.Range("F15").Value = 70
$wsExpo.Range("F22").Value = 1509

# --- 全部类型 (All types, combined listing) sheet -------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value   = 242
$wsAll.Range("F7").Value   = 9995
$wsAll.Range("F9").Value   = 331
$wsAll.Range("F10").Value  = 1244
$wsAll.Range("F11").Value  = 5340
$wsAll.Range("F15").Value  = 121
$wsAll.Range("F16").Value  = 70
$wsAll.Range("F23").Value  = 1509
